$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (it will be re-created at the
#    new last-edit location, the paragraph we are about to insert).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Insert a new list paragraph right after paragraph 7
#    ("Run php artisan key:generate command in CMD"), containing:
#    "Run **npm** **run dev** command in **CMD**" plus the _GoBack bookmark
#    at its end. We use InsertXML so we get full control over run
#    boundaries, bold formatting and the spell-check proof marks exactly
#    like the target markup.
$p7 = $d.Paragraphs.Item(7)
$insertionPoint = $d.Range($p7.Range.End, $p7.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">Run </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>npm</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> run dev </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">command in </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>CMD</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p>' +
  '<w:sectPr></w:sectPr>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml)

# InsertXML needed a trailing paragraph mark to force the split from the
# following paragraph; that left behind one extra empty list paragraph.
# Remove it (together with its own paragraph mark) now that the split is
# in place.
$emptyPara = $d.Paragraphs.Item(9)
$emptyPara.Range.Delete()
